$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Model sheet: add HR_DUTY_ROSTER row (row 4)
# ---------------------------------------------------------------------------
$wsModel = $wb.Worksheets.Item("Model")
$wsModel.Range("A4").Value = "HR_DUTY_ROSTER"
$wsModel.Range("B4").Value = "Yes"
$wsModel.Range("C4").Value = "Yes"
$wsModel.Range("D4").Value = "No"
$wsModel.Range("E4").Value = "Yes"
$wsModel.Range("F4").Value = "Yes"
$wsModel.Range("G4").Value = "Yes"
$wsModel.Range("H4").Value = "Yes"
$wsModel.Range("A3:H3").Copy()
$wsModel.Range("A4:H4").PasteSpecial(-4122)
$wsModel.Range("G9").Select()

# ---------------------------------------------------------------------------
# Controller sheet: add DutyRosterController row (row 4)
# ---------------------------------------------------------------------------
$wsController = $wb.Worksheets.Item("Controller")
$wsController.Range("A4").Value = "DutyRosterController"
$wsController.Range("B4").Value = "Yes"
$wsController.Range("C4").Value = "Yes"
$wsController.Range("D4").Value = "Yes"
$wsController.Range("E4").Value = "Yes"
$wsController.Range("F4").Value = "Yes"
$wsController.Range("G4").Value = "Yes"
$wsController.Range("H4").Value = "Yes"
$wsController.Range("I4").Value = "Yes"
$wsController.Range("J4").Value = "Yes"
$wsController.Range("K4").Value = "No"
$wsController.Range("L4").Value = "No"
$wsController.Range("A3:L3").Copy()
$wsController.Range("A4:L4").PasteSpecial(-4122)
$wsController.Range("B4").Select()

# ---------------------------------------------------------------------------
# Views sheet: add DutyRoster row (row 4)
# ---------------------------------------------------------------------------
$wsViews = $wb.Worksheets.Item("Views")
$wsViews.Range("A4").Value = "DutyRoster"
$wsViews.Range("B4").Value = "Yes"
$wsViews.Range("C4").Value = "Yes"
$wsViews.Range("D4").Value = "Yes"
$wsViews.Range("E4").Value = "Yes"
$wsViews.Range("F4").Value = "No"
$wsViews.Range("G4").Value = "No"
$wsViews.Range("H4").Value = "No"
$wsViews.Range("I4").Value = "Yes"
$wsViews.Range("J4").Value = "Yes"
$wsViews.Range("K4").Value = "Yes"
$wsViews.Range("L4").Value = "No"
$wsViews.Range("M4").Value = "No"
$wsViews.Range("N4").Value = "No"
$wsViews.Range("A3:N3").Copy()
$wsViews.Range("A4:N4").PasteSpecial(-4122)
$wsViews.Range("A3").Select()

$wsViews.Activate()
